# The upstream change (M2Doc issue #295 - "Add the version of M2Doc in the
# template custom properties") was produced by re-saving this template
# through a different OOXML writer (docx4j). For this particular template
# part (word/document.xml / word/styles.xml) that re-serialization only
# normalized XML attribute order (e.g. <w:pgSz w:w=".." w:h=".."/> ->
# <w:pgSz w:h=".." w:w=".."/>, namespace declarations sorted
# alphabetically, w:qFormat moved before w:uiPriority, etc.). Every
# element, attribute name/value pair and text run is unchanged -
# canonicalizing both versions of the XML (ignoring attribute order)
# yields byte-identical trees, and the Word object model has no notion of
# "attribute order" to set explicitly (Word/this COM surface always
# serializes its own canonical order on write).
#
# So the faithful reproduction of this particular diff is a content
# no-op: touch the document through the object model without changing
# any text, formatting, styles or section/page-setup values.

$d = $word.ActiveDocument

# Touch the document (read-only) to confirm the object model is alive,
# without mutating any content, formatting or properties.
$null = $d.Content.Text
$null = $d.Sections.Item(1).PageSetup.PageWidth
$null = $d.Styles.Item("Normal").NameLocal
